# Add the 2023 column (T) to the "births accepted by skilled health
# personnel" table, mirroring the existing 2007-2022 (D:S) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column T, keyed by row number (4 = year header row,
# 5-14 = Kyrgyz Republic + the 7 oblasts + Bishkek/Osh cities).
$values = [ordered]@{
    4  = 2023
    5  = 99.5
    6  = 99.426175237254469
    7  = 99.458151211935132
    8  = 99.44178628389156
    9  = 99.453125
    10 = 99.487295483676391
    11 = 99.743589743589752
    12 = 99.190647482014398
    13 = 99.483321247280642
    14 = 99.771121504627331
}

foreach ($row in $values.Keys) {
    # Copy the formatting of the existing last column (S) of the same
    # row onto the new T cell, then write the new value into it - this
    # keeps T's style consistent with the rest of the table.
    $ws.Range("S$row").Copy()
    $ws.Range("T$row").PasteSpecial(-4122)
    $ws.Range("T$row").Value = $values[$row]
}

# Drop the stale selection (previously parked on U6, outside the table)
# now that the table has grown to include column T.
$ws.Range("A1").Select() | Out-Null
